# Update "想去人数" (F column) counts on multiple sheets, reflecting a
# refreshed data scrape (gh-pages output regeneration at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 11
$ws.Range("F5").Value = 74
$ws.Range("F6").Value = 881
$ws.Range("F7").Value = 54
$ws.Range("F8").Value = 6980
$ws.Range("F11").Value = 148
$ws.Range("F12").Value = 6491
$ws.Range("F15").Value = 4434
$ws.Range("F17").Value = 52
$ws.Range("F18").Value = 4463
$ws.Range("F19").Value = 7
$ws.Range("F20").Value = 240
$ws.Range("F21").Value = 248
$ws.Range("F23").Value = 231
$ws.Range("F25").Value = 172
$ws.Range("F29").Value = 8036
$ws.Range("F31").Value = 1386
$ws.Range("F32").Value = 57
$ws.Range("F33").Value = 692
$ws.Range("F37").Value = 69
$ws.Range("F38").Value = 1626
$ws.Range("F39").Value = 210
$ws.Range("F40").Value = 938
$ws.Range("F41").Value = 41
$ws.Range("F42").Value = 4081
$ws.Range("F43").Value = 347
$ws.Range("F44").Value = 25
$ws.Range("F46").Value = 45
$ws.Range("F47").Value = 836
$ws.Range("F48").Value = 1109
$ws.Range("F49").Value = 13

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 16
$ws.Range("F19").Value = 873

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 237

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 237
$ws.Range("F6").Value = 11
$ws.Range("F8").Value = 74
$ws.Range("F9").Value = 16
$ws.Range("F10").Value = 881
$ws.Range("F11").Value = 54
$ws.Range("F12").Value = 6980
$ws.Range("F15").Value = 148
$ws.Range("F16").Value = 6491
$ws.Range("F19").Value = 4434
$ws.Range("F21").Value = 4463
$ws.Range("F22").Value = 7
$ws.Range("F23").Value = 240
$ws.Range("F24").Value = 248
$ws.Range("F26").Value = 231
$ws.Range("F29").Value = 8036
$ws.Range("F31").Value = 1386
$ws.Range("F32").Value = 57
$ws.Range("F33").Value = 692
$ws.Range("F37").Value = 69
$ws.Range("F38").Value = 1626
$ws.Range("F39").Value = 210
$ws.Range("F40").Value = 938
$ws.Range("F41").Value = 41
$ws.Range("F42").Value = 4081
$ws.Range("F43").Value = 347
$ws.Range("F44").Value = 25
$ws.Range("F46").Value = 836
$ws.Range("F47").Value = 1109
$ws.Range("F49").Value = 13
